$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 66, shifting rows 66-73 down to 67-74.
$ws.Rows.Item(66).Insert()

# Fill in the new row 66 with data copied from row 65 (now shifted below),
# but with the date changed to 44783.
$ws.Cells.Item(66, 1).Value = 9
$ws.Cells.Item(66, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(66, 3).Value = "Metropolitana"
$ws.Cells.Item(66, 4).Value = 44783
$ws.Cells.Item(66, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(66, 5).Value = 13
$ws.Cells.Item(66, 6).Value = 100112029
$ws.Cells.Item(66, 7).Value = "Orégano"
$ws.Cells.Item(66, 8).Value = "Sin especificar"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 16
$ws.Cells.Item(66, 11).Value = 20000
$ws.Cells.Item(66, 12).Value = 20000
$ws.Cells.Item(66, 13).Value = 20000
$ws.Cells.Item(66, 14).Value = "$/docena de atados"
$ws.Cells.Item(66, 15).Value = "Región Metropolitana"
$ws.Cells.Item(66, 16).Value = 6667
$ws.Cells.Item(66, 17).Value = 3
$ws.Cells.Item(66, 18).Value = "Hortaliza"
